$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44648
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 6750
$ws.Range("P2").Value = 112
$ws.Range("D3").Value = 44967
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 4500
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 4850
$ws.Range("N3").Value = "$/caja 90 unidades"
$ws.Range("P3").Value = 54
$ws.Range("Q3").Value = 90
$ws.Range("D4").Value = 44589
$ws.Range("J4").Value = 110
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = 5500
$ws.Range("P4").Value = 92
$ws.Range("D5").Value = 44421
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8500
$ws.Range("P5").Value = 142
$ws.Range("D6").Value = 44764
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 7500
$ws.Range("P6").Value = 125
$ws.Range("D7").Value = 44669
$ws.Range("J7").Value = 130
$ws.Range("K7").Value = 4500
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = 4750
$ws.Range("P7").Value = 79
$ws.Range("D8").Value = 45044
$ws.Range("J8").Value = 190
$ws.Range("K8").Value = 4000
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = 4526
$ws.Range("P8").Value = 75
$ws.Range("D9").Value = 44657
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 5500
$ws.Range("M9").Value = 5250
$ws.Range("P9").Value = 88
$ws.Range("D10").Value = 44760
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7500
$ws.Range("M10").Value = 7250
$ws.Range("P10").Value = 121
$ws.Range("D11").Value = 44676
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 4500
$ws.Range("M11").Value = 4250
$ws.Range("P11").Value = 71
$ws.Range("D12").Value = 44362
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 8500
$ws.Range("P12").Value = 142
$ws.Range("D13").Value = 44627
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4500
$ws.Range("M13").Value = 4250
$ws.Range("P13").Value = 71
$ws.Range("D14").Value = 44242
$ws.Range("J14").Value = 160
$ws.Range("L14").Value = 5500
$ws.Range("M14").Value = 5250
$ws.Range("P14").Value = 88
$ws.Range("D15").Value = 44494
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 6000
$ws.Range("M15").Value = 5500
$ws.Range("P15").Value = 92
$ws.Range("D16").Value = 44827
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 6500
$ws.Range("P16").Value = 108
$ws.Range("D17").Value = 44963
$ws.Range("J17").Value = 130
$ws.Range("K17").Value = 4000
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = 4250
$ws.Range("P17").Value = 71
$ws.Range("D18").Value = 44400
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 9000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 9500
$ws.Range("N18").Value = "$/caja 60 unidades"
$ws.Range("P18").Value = 158
$ws.Range("Q18").Value = 60
$ws.Range("D19").Value = 45079
$ws.Range("J19").Value = 130
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 4462
$ws.Range("P19").Value = 74
$ws.Range("D20").Value = 44935
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 6500
$ws.Range("P20").Value = 108
$ws.Range("D21").Value = 44281
$ws.Range("K21").Value = 5500
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = 5750
$ws.Range("P21").Value = 96
$ws.Range("D22").Value = 44785
$ws.Range("J22").Value = 130
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 7500
$ws.Range("P22").Value = 125
$ws.Range("D23").Value = 44740
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 6500
$ws.Range("P23").Value = 108
$ws.Range("D24").Value = 44382
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 8000
$ws.Range("M24").Value = 7438
$ws.Range("P24").Value = 124
$ws.Range("D25").Value = 44603
$ws.Range("J25").Value = 140
$ws.Range("K25").Value = 5500
$ws.Range("M25").Value = 5750
$ws.Range("P25").Value = 96
